$d = $word.ActiveDocument
$paras = $d.Paragraphs

# Locate the two target paragraphs by their distinctive text so the edit
# only touches the "Methods or software..." / "Description of protocol..."
# sentences and not the look-alike "[targetaudience]" placeholder that
# still lives in the summary table further down the document.
$targetAudiencePara = $null
$toolsPara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($t.Contains("Methods or software needed to access and use data:") -and $t.Contains("[targetaudience]")) {
        $targetAudiencePara = $p
    } elseif ($t.Contains("Description of protocol to access restricted data:") -and $t.Contains("[tools]")) {
        $toolsPara = $p
    }
}

# 1) "[targetaudience]" -> "[tools]" in the "Methods or software..." sentence.
if ($targetAudiencePara -ne $null) {
    $rng = $targetAudiencePara.Range
    $rng.Find.Execute("[targetaudience]", $false, $false, $false, $false, $false,
                       $true, 1, $false, "[tools]", 2)
}

# 2) Drop the trailing " [tools]" run from the
#    "Description of protocol to access restricted data:" sentence, leaving
#    just the ":" run behind.
if ($toolsPara -ne $null) {
    $rng2 = $toolsPara.Range
    $rng2.Find.Execute(" [tools]", $false, $false, $false, $false, $false,
                        $true, 1, $false, "", 2)
}
